# Regenerate save_data column G ("K") values: use K instead of Strike#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 4
    4  = 2
    5  = 6
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 5
    11 = 1
    12 = 2
    13 = 2
    14 = 0
    15 = 5
    16 = 2
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
